$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8992.429
$ws.Range("I32").Value = 8986.75
$ws.Range("K32").Value = 8986.75
$ws.Range("M32").Value = -8660.75
$ws.Range("H40").Value = 1790.6666
$ws.Range("I40").Value = 1166.4445
$ws.Range("J40").Value = 3663.3333
$ws.Range("K40").Value = 1166.4445
$ws.Range("L40").Value = 3663.3333
$ws.Range("M40").Value = -991.4445000000001
$ws.Range("N40").Value = -4013.3333
$ws.Range("H100").Value = 2119
$ws.Range("I100").Value = 2143.3333
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 2143.3333
$ws.Range("L100").Value = 1900
$ws.Range("M100").Value = -1602.3333
$ws.Range("N100").Value = -2982

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 866.5
$ws.Range("I97").Value = 204.57143
$ws.Range("K97").Value = 204.57143
$ws.Range("M97").Value = 291.42857
$ws.Range("H102").Value = 1774.75
$ws.Range("I102").Value = 1799.7273
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1799.7273
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -177.7273
$ws.Range("N102").Value = -4744
$ws.Range("H122").Value = 1362.3334
$ws.Range("I122").Value = 1411.871
$ws.Range("K122").Value = 4235.613
$ws.Range("M122").Value = -1785.613
$ws.Range("H132").Value = 5819.15
$ws.Range("I132").Value = 6022.5884
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 18067.7652
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = -15537.7652
$ws.Range("N132").Value = -19059.0005

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2386.6155
$ws.Range("I16").Value = 2355.875
$ws.Range("K16").Value = 2355.875
$ws.Range("M16").Value = -2068.875
$ws.Range("H99").Value = 3281.2856
$ws.Range("I99").Value = 2994.8333
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 2994.8333
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -1496.8333
$ws.Range("N99").Value = -7996
$ws.Range("H107").Value = 1263.6923
$ws.Range("I107").Value = 557.1111
$ws.Range("J107").Value = 2853.5
$ws.Range("K107").Value = 557.1111
$ws.Range("L107").Value = 2853.5
$ws.Range("M107").Value = 1362.8889
$ws.Range("N107").Value = -6693.5
$ws.Range("H113").Value = 2386.6155
$ws.Range("I113").Value = 2355.875
$ws.Range("K113").Value = 2355.875
$ws.Range("M113").Value = -185.875
$ws.Range("H126").Value = 3281.2856
$ws.Range("I126").Value = 2994.8333
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 8984.499899999999
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -6514.499899999999
$ws.Range("N126").Value = -19940
$ws.Range("H134").Value = 3487.3076
$ws.Range("I134").Value = 3487.3076
$ws.Range("K134").Value = 10461.9228
$ws.Range("M134").Value = -7926.9228

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 221.5
$ws.Range("I92").Value = 221.5
$ws.Range("K92").Value = 664.5
$ws.Range("M92").Value = 583.5
$ws.Range("H129").Value = 2513.4443
$ws.Range("J129").Value = 2513.4443
$ws.Range("L129").Value = 7540.3329
$ws.Range("N129").Value = -17540.3329
$ws.Range("H139").Value = 5498.143
$ws.Range("I139").Value = 5264.5
$ws.Range("K139").Value = 15793.5
$ws.Range("M139").Value = -10653.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17925.908
$ws.Range("J70").Value = 6747
$ws.Range("L70").Value = 6747
$ws.Range("N70").Value = -7287
$ws.Range("H73").Value = 17925.908
$ws.Range("J73").Value = 6747
$ws.Range("L73").Value = 6747
$ws.Range("N73").Value = -8619
$ws.Range("H80").Value = 4333
$ws.Range("I80").Value = 4333
$ws.Range("K80").Value = 4333
$ws.Range("M80").Value = -3335
$ws.Range("H83").Value = 4333
$ws.Range("I83").Value = 4333
$ws.Range("K83").Value = 21665
$ws.Range("M83").Value = -16673
$ws.Range("H99").Value = 110393.8
$ws.Range("I99").Value = 8242.25
$ws.Range("K99").Value = 8242.25
$ws.Range("M99").Value = -5996.25
$ws.Range("H107").Value = 5298.6665
$ws.Range("I107").Value = 9000
$ws.Range("K107").Value = 9000
$ws.Range("M107").Value = -7080
$ws.Range("H122").Value = 3939.25
$ws.Range("I122").Value = 3086.5
$ws.Range("K122").Value = 9259.5
$ws.Range("M122").Value = -6809.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H40").Value = 3154.3333
$ws.Range("I40").Value = 3154.3333
$ws.Range("K40").Value = 3154.3333
$ws.Range("M40").Value = -3018.3333
$ws.Range("H46").Value = 2245
$ws.Range("I46").Value = 1194
$ws.Range("J46").Value = 7500
$ws.Range("K46").Value = 1194
$ws.Range("L46").Value = 7500
$ws.Range("M46").Value = -1006
$ws.Range("N46").Value = -7876
$ws.Range("H61").Value = 1795.5
$ws.Range("J61").Value = 1399
$ws.Range("L61").Value = 1399
$ws.Range("N61").Value = -1803
$ws.Range("H113").Value = 1795.5
$ws.Range("J113").Value = 1399
$ws.Range("L113").Value = 1399
$ws.Range("N113").Value = -5739
$ws.Range("H132").Value = 4781.222
$ws.Range("I132").Value = 4504.4287
$ws.Range("K132").Value = 13513.2861
$ws.Range("M132").Value = -10983.2861
$ws.Range("H136").Value = 3021.647
$ws.Range("I136").Value = 1513
$ws.Range("J136").Value = 7924.75
$ws.Range("K136").Value = 4539
$ws.Range("L136").Value = 23774.25
$ws.Range("M136").Value = -1989
$ws.Range("N136").Value = -28874.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 386
$ws.Range("I107").Value = 386
$ws.Range("K107").Value = 1158
$ws.Range("M107").Value = 762
$ws.Range("H113").Value = 233.33333
$ws.Range("I113").Value = 275
$ws.Range("K113").Value = 825
$ws.Range("M113").Value = 1345
$ws.Range("H122").Value = 4013.5
$ws.Range("I122").Value = 3790.9167
$ws.Range("K122").Value = 11372.7501
$ws.Range("M122").Value = -8922.750100000001
$ws.Range("H136").Value = 8215.565000000001
$ws.Range("I136").Value = 5676.3335
$ws.Range("K136").Value = 17029.0005
$ws.Range("M136").Value = -14479.0005
